$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: add a labeled_h5 value (D7) that was missing before ---
$ws.Range("D7").Value = "CollectedData_wi.h5"
$ws.Range("D7").WrapText = $true
$ws.Range("D7").VerticalAlignment = -4160

# --- Row 8 (new): stage=initial_labeling, labeling_path=test03 folder, comment "add spine" ---
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "initial_labeling"
$ws.Range("C8").Value = "W:\wataru\dlc_data\homecage_test03-wi-2022-06-08\labeled-data\rpicam-01_1806_20210722_212134"
$ws.Range("J8").Value = "add spine"

# --- Row 9 (new): stage=initial_labeling, comment "spine1,spine2,tail1,tail2", labeling_path=test04 folder ---
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "initial_labeling"
$ws.Range("J9").Value = "spine1,spine2,tail1,tail2"
$ws.Range("C9").Value = "W:\wataru\dlc_data\homecage_test04-wi-2022-06-08\labeled-data\rpicam-01_1806_20210722_212134"

# --- fill the remaining (reused) values for rows 8 and 9 ---
$ws.Range("D8").Value = "CollectedData_wi.h5"
$ws.Range("G8").Value = "W:\wataru\homecage_videos\black_mice\13_pair\analyzed_videos"
$ws.Range("H8").Value = "rpicam-01_1806_20210722_212134.mp4"
$ws.Range("K8").Value = ""

$ws.Range("D9").Value = "CollectedData_wi.h5"
$ws.Range("G9").Value = "W:\wataru\homecage_videos\black_mice\13_pair\analyzed_videos"
$ws.Range("H9").Value = "rpicam-01_1806_20210722_212134.mp4"

# --- formatting: wrap text, top vertical alignment, row height for new rows ---
$ws.Range("A8:D8").WrapText = $true
$ws.Range("A8:D8").VerticalAlignment = -4160
$ws.Range("G8:H8").WrapText = $true
$ws.Range("G8:H8").VerticalAlignment = -4160
$ws.Range("J8:K8").WrapText = $true
$ws.Range("J8:K8").VerticalAlignment = -4160
$ws.Rows.Item(8).RowHeight = 87

$ws.Range("A9:D9").WrapText = $true
$ws.Range("A9:D9").VerticalAlignment = -4160
$ws.Range("G9:H9").WrapText = $true
$ws.Range("G9:H9").VerticalAlignment = -4160
$ws.Range("J9").WrapText = $false
$ws.Range("J9").VerticalAlignment = -4160
$ws.Rows.Item(9).RowHeight = 87

# --- view: freeze header row and scroll so row 7 is at the top, selection on J18 ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("J18").Select()
